$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 155/156, pushing existing rows 155-242 down to 157-244.
$ws.Rows("155:156").Insert()

# --- New row 155 ---
$ws.Cells.Item(155, 1).Value2 = 5
$ws.Cells.Item(155, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(155, 3).Value2 = "Maule"
$ws.Cells.Item(155, 4).Value2 = 44488
$ws.Cells.Item(155, 5).Value2 = 7
$ws.Cells.Item(155, 6).Value2 = 100112043
$ws.Cells.Item(155, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(155, 8).Value2 = "Sin especificar"
$ws.Cells.Item(155, 9).Value2 = "Primera"
$ws.Cells.Item(155, 10).Value2 = 300
$ws.Cells.Item(155, 11).Value2 = 9000
$ws.Cells.Item(155, 12).Value2 = 9000
$ws.Cells.Item(155, 13).Value2 = 9000
$ws.Cells.Item(155, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(155, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value2 = 150
$ws.Cells.Item(155, 17).Value2 = 60
$ws.Cells.Item(155, 18).Value2 = "Hortaliza"

# --- New row 156 ---
$ws.Cells.Item(156, 1).Value2 = 5
$ws.Cells.Item(156, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(156, 3).Value2 = "Maule"
$ws.Cells.Item(156, 4).Value2 = 44488
$ws.Cells.Item(156, 5).Value2 = 7
$ws.Cells.Item(156, 6).Value2 = 100112043
$ws.Cells.Item(156, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(156, 8).Value2 = "Sin especificar"
$ws.Cells.Item(156, 9).Value2 = "Primera"
$ws.Cells.Item(156, 10).Value2 = 400
$ws.Cells.Item(156, 11).Value2 = 12000
$ws.Cells.Item(156, 12).Value2 = 12000
$ws.Cells.Item(156, 13).Value2 = 12000
$ws.Cells.Item(156, 14).Value2 = "$/caja 80 unidades"
$ws.Cells.Item(156, 15).Value2 = "Región del Maule"
$ws.Cells.Item(156, 16).Value2 = 150
$ws.Cells.Item(156, 17).Value2 = 80
$ws.Cells.Item(156, 18).Value2 = "Hortaliza"
